$d = $word.ActiveDocument

# Bold the glossary "term" that precedes each " <en-dash> definition"
# sentence in the "Highlighted vocabulary" list, splitting each paragraph's
# single run into a bold term run + a regular-weight definition run.
$glossary = @(
    @("Cholesterol", " – it is a waxy, fat-like substance that is found in all the cells in your body."),
    @("Hormones", " – help control how cells and organs do their work in your body e.g., insulin."),
    @("Reduces", " – to make smaller/less."),
    @("Feel-good", " – a feeling of happiness."),
    @("Increases", " – to make bigger/more."),
    @("Mental health", " – your emotional state."),
    @("Mental awareness", " – to make us aware of the importance of being happy and feel good and not being depressed and sad."),
    @("Stroking your pet", " – rubbing your hand over your pet's fur.")
)

foreach ($pair in $glossary) {
    $term = $pair[0]
    $rest = $pair[1]

    $rng = $d.Content
    $found = $rng.Find.Execute($term + $rest, $true, $false, $false, $false, $false, $true, 1, $false, "", 0)
    if (-not $found) {
        continue
    }

    $start = $rng.Start
    $boldLen = $term.Length + 1
    $boldRng = $d.Range($start, $start + $boldLen)
    $boldRng.Font.Bold = 1
}
